$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 149, shifting existing rows 149:200 down to 150:201
$ws.Rows.Item(149).Insert()

# Populate the newly inserted row 149 with its data (matching the style/pattern of surrounding rows)
$ws.Cells.Item(149, 1).Value = 5
$ws.Cells.Item(149, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(149, 3).Value = "Maule"
$ws.Cells.Item(149, 4).Value = 44468
$ws.Cells.Item(149, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(149, 5).Value = 7
$ws.Cells.Item(149, 6).Value = 100112032
$ws.Cells.Item(149, 7).Value = "Zapallo italiano"
$ws.Cells.Item(149, 8).Value = "Sin especificar"
$ws.Cells.Item(149, 9).Value = "Primera"
$ws.Cells.Item(149, 10).Value = 200
$ws.Cells.Item(149, 11).Value = 18000
$ws.Cells.Item(149, 12).Value = 18000
$ws.Cells.Item(149, 13).Value = 18000
$ws.Cells.Item(149, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(149, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(149, 16).Value = 360
$ws.Cells.Item(149, 17).Value = 50
$ws.Cells.Item(149, 18).Value = "Hortaliza"
